$wb = $excel.ActiveWorkbook

# --- Monsters sheet (sheet9.xml) ---
$wsMonsters = $wb.Worksheets.Item("Monsters")

# Row 4: new monster entry in left table (A:C)
$wsMonsters.Range("A4").Value = 59
$wsMonsters.Range("B4").Value = "Untoter Krieger"
$wsMonsters.Range("C4").Value = "Monster in Manyeyes' castle"

# Row 5: new monster entry in left table (A:C) + right table (H:I)
$wsMonsters.Range("A5").Value = 60
$wsMonsters.Range("B5").Value = "Untoter Magier"
$wsMonsters.Range("C5").Value = "Monster in Manyeyes' castle"
$wsMonsters.Range("H5").Value = 90
$wsMonsters.Range("I5").Value = "2x Untoter Krieger"

# Row 6: right table (H:I) only
$wsMonsters.Range("H6").Value = 91
$wsMonsters.Range("I6").Value = "2x Untoter Krieger, 1x Untoter Magier"

# Row 7: right table (H:I) only
$wsMonsters.Range("H7").Value = 92
$wsMonsters.Range("I7").Value = "3x Untoter Krieger, 2x Untoter Magier"

# Column C autofit (becomes bestFit with wider width due to new longer text)
$wsMonsters.Columns.Item(3).AutoFit()

# --- GlobalVars sheet (sheet3.xml) loses the tabSelected flag, gets new selection ---
$wsGlobalVars = $wb.Worksheets.Item("GlobalVars")
$wsGlobalVars.Range("D16").Select()

# --- Monsters sheet becomes the active tab, with new selection ---
$wsMonsters.Activate()
$wsMonsters.Range("I8").Select()
